$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update geothermal_counts (column Z) from "no" to "yes" for the rows where
# additional next-capacity functions now report a geothermal build, so the
# results replicate the previous run.
$ws.Range("Z2").Value = "yes"
$ws.Range("Z4").Value = "yes"
$ws.Range("Z5").Value = "yes"
$ws.Range("Z6").Value = "yes"
$ws.Range("Z7").Value = "yes"
$ws.Range("Z8").Value = "yes"
$ws.Range("Z9").Value = "yes"
$ws.Range("Z10").Value = "yes"
$ws.Range("Z11").Value = "yes"
$ws.Range("Z12").Value = "yes"
$ws.Range("Z13").Value = "yes"

# El Paso (row 13) now gets a non-zero wind capacity factor from the
# additional next capacity function, matching the previous simulation.
$ws.Range("T13").Value = 0.34548755599999997
$ws.Range("V13").Value = 0.34548755599999997
$ws.Range("Y13").Value = 0

# Restore the active selection left by the author after the edit.
$ws.Range("A7").Select()
